$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 27
$ws.Range("I2").ClearContents()
$ws.Range("J2").Value = "new"
$ws.Range("K2").Value = "f"
$ws.Range("L2").Value = "stimuli/img_psgf7.png"
$ws.Range("M2").Value = 26
$ws.Range("N2").Value = 11.66666666666667
$ws.Range("O2").Value = 18.83333333333333
$ws.Range("P2").Value = 36
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("F3").Value = 28
$ws.Range("L3").Value = "stimuli/img_t4hvr.png"
$ws.Range("M3").Value = 61.69230769230769
$ws.Range("N3").Value = 39.76923076923077
$ws.Range("O3").Value = 50.73076923076923
$ws.Range("P3").Value = 39
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3
$ws.Range("F4").Value = 29
$ws.Range("L4").Value = "stimuli/img_yteqw.png"
$ws.Range("M4").Value = 66.83783783783784
$ws.Range("N4").Value = 43.78378378378378
$ws.Range("O4").Value = 55.31081081081081
$ws.Range("P4").Value = 37
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = 4
$ws.Range("S4").Value = 4
$ws.Range("F5").Value = 30
$ws.Range("L5").Value = "stimuli/img_aweye.png"
$ws.Range("M5").Value = 53.42105263157895
$ws.Range("N5").Value = 31.84210526315789
$ws.Range("O5").Value = 42.63157894736842
$ws.Range("P5").Value = 38
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2
$ws.Range("F6").Value = 31
$ws.Range("I6").ClearContents()
$ws.Range("J6").Value = "new"
$ws.Range("K6").Value = "f"
$ws.Range("L6").Value = "stimuli/img_sltwe.png"
$ws.Range("M6").Value = 72.02500000000001
$ws.Range("N6").Value = 46.875
$ws.Range("O6").Value = 59.45
$ws.Range("P6").Value = 40
$ws.Range("Q6").Value = 5
$ws.Range("R6").Value = 5
$ws.Range("S6").Value = 5
$ws.Range("F7").Value = 32
$ws.Range("I7").ClearContents()
$ws.Range("J7").Value = "new"
$ws.Range("K7").Value = "f"
$ws.Range("L7").Value = "stimuli/img_ozxpp.png"
$ws.Range("M7").Value = 26.26470588235294
$ws.Range("N7").Value = 11.47058823529412
$ws.Range("O7").Value = 18.86764705882353
$ws.Range("P7").Value = 34
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("F8").Value = 33
$ws.Range("L8").Value = "stimuli/img_kzg3h.png"
$ws.Range("M8").Value = 77.02777777777777
$ws.Range("N8").Value = 56.22222222222222
$ws.Range("O8").Value = 66.625
$ws.Range("P8").Value = 36
$ws.Range("Q8").Value = 7
$ws.Range("R8").Value = 7
$ws.Range("S8").Value = 7
$ws.Range("F9").Value = 34
$ws.Range("I9").Value = "target"
$ws.Range("J9").Value = "old"
$ws.Range("K9").Value = "j"
$ws.Range("L9").Value = "stimuli/img_ic3os.png"
$ws.Range("M9").Value = 84.79069767441861
$ws.Range("N9").Value = 66.16279069767442
$ws.Range("O9").Value = 75.47674418604652
$ws.Range("P9").Value = 43
$ws.Range("Q9").Value = 9
$ws.Range("R9").Value = 9
$ws.Range("S9").Value = 9
$ws.Range("F10").Value = 35
$ws.Range("F11").Value = 36
$ws.Range("L11").Value = "stimuli/img_72fmj.png"
$ws.Range("M11").Value = 53.87179487179487
$ws.Range("N11").Value = 36.02564102564103
$ws.Range("O11").Value = 44.94871794871795
$ws.Range("P11").Value = 39
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = 3
$ws.Range("S11").Value = 3
$ws.Range("F12").Value = 37
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("J12").Value = "catch"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/catch_17.jpg"
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("O12").ClearContents()
$ws.Range("P12").ClearContents()
$ws.Range("Q12").ClearContents()
$ws.Range("R12").ClearContents()
$ws.Range("S12").ClearContents()
$ws.Range("F13").Value = 38
$ws.Range("L13").Value = "stimuli/img_gbypq.png"
$ws.Range("M13").Value = 76.27500000000001
$ws.Range("N13").Value = 51.925
$ws.Range("O13").Value = 64.09999999999999
$ws.Range("P13").Value = 40
$ws.Range("Q13").Value = 6
$ws.Range("R13").Value = 6
$ws.Range("S13").Value = 6
$ws.Range("F14").Value = 39
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = "new"
$ws.Range("K14").Value = "f"
$ws.Range("L14").Value = "stimuli/img_4wq98.png"
$ws.Range("M14").Value = 78.48387096774194
$ws.Range("N14").Value = 58.12903225806452
$ws.Range("O14").Value = 68.30645161290323
$ws.Range("P14").Value = 31
$ws.Range("Q14").Value = 7
$ws.Range("R14").Value = 7
$ws.Range("S14").Value = 7
$ws.Range("F15").Value = 40
$ws.Range("L15").Value = "stimuli/img_jge7p.png"
$ws.Range("M15").Value = 90.42424242424242
$ws.Range("N15").Value = 75.63636363636364
$ws.Range("O15").Value = 83.03030303030303
$ws.Range("P15").Value = 33
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = 10
$ws.Range("S15").Value = 10
$ws.Range("F16").Value = 41
$ws.Range("I16").Value = "target"
$ws.Range("J16").Value = "old"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_ose78.png"
$ws.Range("M16").Value = 80.19444444444444
$ws.Range("N16").Value = 60.25
$ws.Range("O16").Value = 70.22222222222223
$ws.Range("Q16").Value = 8
$ws.Range("R16").Value = 7
$ws.Range("S16").Value = 7
$ws.Range("F17").Value = 42
$ws.Range("H17").Value = "bedrooms"
$ws.Range("I17").Value = "target"
$ws.Range("J17").Value = "old"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_fnu4h.png"
$ws.Range("M17").Value = 85.87179487179488
$ws.Range("N17").Value = 70.71794871794872
$ws.Range("O17").Value = 78.2948717948718
$ws.Range("P17").Value = 39
$ws.Range("Q17").Value = 9
$ws.Range("R17").Value = 9
$ws.Range("S17").Value = 9
$ws.Range("F18").Value = 43
$ws.Range("L18").Value = "stimuli/img_2pnl2.png"
$ws.Range("M18").Value = 6.621621621621622
$ws.Range("N18").Value = 7.135135135135135
$ws.Range("O18").Value = 6.878378378378379
$ws.Range("P18").Value = 37
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 1
$ws.Range("S18").Value = 1
$ws.Range("F19").Value = 44
$ws.Range("L19").Value = "stimuli/img_1vq1v.png"
$ws.Range("M19").Value = 69.42857142857143
$ws.Range("N19").Value = 46.59523809523809
$ws.Range("O19").Value = 58.01190476190476
$ws.Range("P19").Value = 42
$ws.Range("Q19").Value = 5
$ws.Range("R19").Value = 5
$ws.Range("S19").Value = 5
$ws.Range("F20").Value = 45
$ws.Range("I20").Value = "target"
$ws.Range("J20").Value = "old"
$ws.Range("K20").Value = "j"
$ws.Range("L20").Value = "stimuli/img_anzgh.png"
$ws.Range("M20").Value = 75.10526315789474
$ws.Range("N20").Value = 55.76315789473684
$ws.Range("O20").Value = 65.43421052631579
$ws.Range("P20").Value = 38
$ws.Range("Q20").Value = 6
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 6
$ws.Range("F21").Value = 46
$ws.Range("L21").Value = "stimuli/img_z3yzz.png"
$ws.Range("M21").Value = 71.71052631578948
$ws.Range("N21").Value = 49.81578947368421
$ws.Range("O21").Value = 60.76315789473685
$ws.Range("P21").Value = 38
$ws.Range("Q21").Value = 5
$ws.Range("R21").Value = 5
$ws.Range("S21").Value = 5
$ws.Range("F22").Value = 47
$ws.Range("I22").Value = "target"
$ws.Range("J22").Value = "old"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_9pfbj.png"
$ws.Range("M22").Value = 91.27272727272727
$ws.Range("N22").Value = 80.09090909090909
$ws.Range("O22").Value = 85.68181818181819
$ws.Range("F23").Value = 48
$ws.Range("L23").Value = "stimuli/img_f4jxo.png"
$ws.Range("M23").Value = 82.91666666666667
$ws.Range("N23").Value = 65.52777777777777
$ws.Range("O23").Value = 74.22222222222223
$ws.Range("P23").Value = 36
$ws.Range("Q23").Value = 8
$ws.Range("R23").Value = 8
$ws.Range("S23").Value = 8
$ws.Range("F24").Value = 49
$ws.Range("L24").Value = "stimuli/img_3bxjb.png"
$ws.Range("M24").Value = 87.28571428571429
$ws.Range("N24").Value = 72.65714285714286
$ws.Range("O24").Value = 79.97142857142858
$ws.Range("P24").Value = 35
$ws.Range("Q24").Value = 10
$ws.Range("R24").Value = 10
$ws.Range("S24").Value = 10
$ws.Range("F25").Value = 50
$ws.Range("L25").Value = "stimuli/img_jivhq.png"
$ws.Range("M25").Value = 37
$ws.Range("N25").Value = 22.26530612244898
$ws.Range("O25").Value = 29.63265306122449
$ws.Range("P25").Value = 49
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 2
$ws.Range("F26").Value = 51
$ws.Range("I26").ClearContents()
$ws.Range("J26").Value = "new"
$ws.Range("K26").Value = "f"
$ws.Range("L26").Value = "stimuli/img_c4uwt.png"
$ws.Range("M26").Value = 44.48387096774194
$ws.Range("N26").Value = 30.06451612903226
$ws.Range("O26").Value = 37.2741935483871
$ws.Range("P26").Value = 31
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = 2
$ws.Range("S26").Value = 2
$ws.Range("F27").Value = 52
$ws.Range("I27").Value = "target"
$ws.Range("J27").Value = "old"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_juob3.png"
$ws.Range("M27").Value = 79.92105263157895
$ws.Range("N27").Value = 59.78947368421053
$ws.Range("O27").Value = 69.85526315789474
$ws.Range("P27").Value = 38
$ws.Range("Q27").Value = 7
$ws.Range("R27").Value = 7
$ws.Range("S27").Value = 7
$ws.Range("F28").Value = 53
$ws.Range("I28").ClearContents()
$ws.Range("J28").Value = "new"
$ws.Range("K28").Value = "f"
$ws.Range("L28").Value = "stimuli/img_zi682.png"
$ws.Range("M28").Value = 84.59999999999999
$ws.Range("N28").Value = 69.52500000000001
$ws.Range("O28").Value = 77.0625
$ws.Range("P28").Value = 40
$ws.Range("Q28").Value = 9
$ws.Range("R28").Value = 9
$ws.Range("S28").Value = 9
$ws.Range("F29").Value = 54
$ws.Range("I29").Value = "target"
$ws.Range("J29").Value = "old"
$ws.Range("K29").Value = "j"
$ws.Range("L29").Value = "stimuli/img_cmyvx.png"
$ws.Range("M29").Value = 64.25
$ws.Range("N29").Value = 40.09375
$ws.Range("O29").Value = 52.171875
$ws.Range("P29").Value = 32
$ws.Range("Q29").Value = 4
$ws.Range("R29").Value = 4
$ws.Range("S29").Value = 4
